$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column B width (23 -> 22)
# NOTE: the host's ColumnWidth setter performs the standard Excel
# character-width -> pixel -> character-width round trip (quantized to
# the active font's "maximum digit width"), so an input of exactly 22
# saves back out as 22.8333. Empirically, 21.15 is the character-width
# input that survives that round trip and lands on a stored width of
# exactly 22.
$ws.Columns.Item(2).ColumnWidth = 21.15

# Update data values in row 2
$ws.Range("B2").Value = 0.000879017197648807
$ws.Range("C2").Value = 0.02820932982320366
$ws.Range("D2").Value = 0.4784197943381507
$ws.Range("E2").Value = 151229

# Update data values in row 3
$ws.Range("B3").Value = 0.001537885844552155
$ws.Range("C3").Value = 0.0276356139230172
$ws.Range("D3").Value = 0.4835471064284624
$ws.Range("E3").Value = 14060
